# Update CONNECTICUT_2023 MCAS sheet:
#  - rename header columns to short machine-friendly names
#  - convert the state/municipality text columns from ALL CAPS to Title Case
#  - drop the trailing sample-size / source / author / date footer rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helpers -----------------------------------------------------------
function IsLetterChar($ch) {
    $up = $ch.ToUpper()
    $lo = $ch.ToLower()
    if ($up.Equals($lo)) {
        return $false
    }
    return $true
}

# Mimic Python's str.title(): capitalize the first letter of every run of
# letters (a "word" = maximal run of letters; anything else, incl. spaces,
# periods and slashes, is a separator) and lowercase the rest of the run.
function PyTitleCase($s) {
    $result = ""
    $prevIsLetter = $false
    for ($i = 0; $i -lt $s.Length; $i++) {
        $ch = $s.Substring($i, 1)
        if (IsLetterChar $ch) {
            if ($prevIsLetter) {
                $result = $result + $ch.ToLower()
            } else {
                $result = $result + $ch.ToUpper()
            }
            $prevIsLetter = $true
        } else {
            $result = $result + $ch
            $prevIsLetter = $false
        }
    }
    return $result
}

# --- 1. header renames ---------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. title-case the state (A) and municipality (B) columns ------------
for ($r = 2; $r -le 521; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    if ($a -ne $null) {
        $ws.Cells.Item($r, 1).Value = PyTitleCase $a
    }
    $b = $ws.Cells.Item($r, 2).Value2
    if ($b -ne $null) {
        $ws.Cells.Item($r, 2).Value = PyTitleCase $b
    }
}

# --- 3. drop the footer rows (blank row 522 + metadata rows 523-527) -----
$ws.Rows("522:527").Delete()
